$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.341.35'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.05%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.76'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.93%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.20'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E5').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4442'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -6.81%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3816'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.33%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.12'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -10.17%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07816'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.76%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.009'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.99%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.28'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.78%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.842.65'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.30%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.828'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.64%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.061'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -7.63%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '85.30'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.27%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001024'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.43%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06511'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.88'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -9.58%  '
$ws.Range('E20').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E21').ClearFormats()
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.440'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.62%  '
$ws.Range('E22').ClearFormats()
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.342.37'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.09%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.74'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -7.23%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.254'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.066.80'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.72%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.28'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.85%  '
$ws.Range('E28').ClearFormats()
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.509'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.42%  '
$ws.Range('E29').ClearFormats()
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.040'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.82%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.22'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09290'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E32').ClearFormats()
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.449'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E33').ClearFormats()
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9237'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.65%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.627'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.217'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.19%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02209'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.29%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05930'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.93%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.195'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.82%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.261'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.42%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5869'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.64%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.24'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.61%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1841'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.50%  '
$ws.Range('E44').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.67%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5614'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.20%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -7.09%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.356'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.914'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.80%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06855'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '107.82'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.70%  '
$ws.Range('E51').ClearFormats()
